$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2187.5715
$ws.Range("I19").Value = 1184.6
$ws.Range("K19").Value = 1184.6
$ws.Range("M19").Value = -1009.6
$ws.Range("H28").Value = 610.55554
$ws.Range("I28").Value = 215.83333
$ws.Range("K28").Value = 215.83333
$ws.Range("M28").Value = 269.16667
$ws.Range("H32").Value = 2898.5454
$ws.Range("J32").Value = 2964.889
$ws.Range("L32").Value = 2964.889
$ws.Range("N32").Value = -3616.889
$ws.Range("H38").Value = 2167.3333
$ws.Range("I38").Value = 1034.6666
$ws.Range("K38").Value = 3103.9998
$ws.Range("M38").Value = -2731.9998
$ws.Range("H40").Value = 2261.125
$ws.Range("J40").Value = 2633
$ws.Range("L40").Value = 2633
$ws.Range("N40").Value = -2983
$ws.Range("H42").Value = 1236.8
$ws.Range("I42").Value = 147
$ws.Range("J42").Value = 1509.25
$ws.Range("K42").Value = 441
$ws.Range("L42").Value = 4527.75
$ws.Range("M42").Value = -211
$ws.Range("N42").Value = -4987.75
$ws.Range("H53").Value = 1759.5
$ws.Range("J53").Value = 95
$ws.Range("L53").Value = 95
$ws.Range("N53").Value = -1369
$ws.Range("H92").Value = 15625227
$ws.Range("I92").Value = 16666842
$ws.Range("J92").Value = 998
$ws.Range("K92").Value = 16666842
$ws.Range("L92").Value = 998
$ws.Range("M92").Value = -16665594
$ws.Range("N92").Value = -3494
$ws.Range("H98").Value = 1552.5938
$ws.Range("I98").Value = 1347.1305
$ws.Range("K98").Value = 1347.1305
$ws.Range("M98").Value = 150.8695
$ws.Range("H106").Value = 2924.5
$ws.Range("I106").Value = 3477.2222
$ws.Range("K106").Value = 3477.2222
$ws.Range("M106").Value = -2846.2222
$ws.Range("H107").Value = 951.6667
$ws.Range("I107").Value = 509.35715
$ws.Range("J107").Value = 2499.75
$ws.Range("K107").Value = 509.35715
$ws.Range("L107").Value = 2499.75
$ws.Range("M107").Value = 1410.64285
$ws.Range("N107").Value = -6339.75
$ws.Range("H122").Value = 1552.5938
$ws.Range("I122").Value = 1347.1305
$ws.Range("K122").Value = 4041.3915
$ws.Range("M122").Value = -1591.3915
$ws.Range("H132").Value = 1003
$ws.Range("I132").Value = 952.78125
$ws.Range("J132").Value = 1324.4
$ws.Range("K132").Value = 2858.34375
$ws.Range("L132").Value = 3973.2
$ws.Range("M132").Value = -328.34375
$ws.Range("N132").Value = -9033.200000000001
$ws.Range("H138").Value = 2339.0635
$ws.Range("I138").Value = 2242.4443
$ws.Range("J138").Value = 2467.889
$ws.Range("K138").Value = 6727.3329
$ws.Range("L138").Value = 7403.667
$ws.Range("M138").Value = -1587.3329
$ws.Range("N138").Value = -17683.667
$ws.Range("H141").Value = 4668940
$ws.Range("I141").Value = 7000771.5
$ws.Range("K141").Value = 21002314.5
$ws.Range("M141").Value = -20997134.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5479.961
$ws.Range("I32").Value = 3553.2537
$ws.Range("J32").Value = 18388.9
$ws.Range("K32").Value = 3553.2537
$ws.Range("L32").Value = 18388.9
$ws.Range("M32").Value = -3266.2537
$ws.Range("N32").Value = -18962.9
$ws.Range("H97").Value = 1379.826
$ws.Range("I97").Value = 1383.125
$ws.Range("J97").Value = 1372.2858
$ws.Range("K97").Value = 1383.125
$ws.Range("L97").Value = 1372.2858
$ws.Range("M97").Value = -887.125
$ws.Range("N97").Value = -2364.2858
$ws.Range("H122").Value = 2348.25
$ws.Range("I122").Value = 1284.875
$ws.Range("K122").Value = 3854.625
$ws.Range("M122").Value = -1404.625
$ws.Range("H130").Value = 53957.25
$ws.Range("J130").Value = 53957.25
$ws.Range("L130").Value = 53957.25
$ws.Range("N130").Value = -63997.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2795.6
$ws.Range("I20").Value = 2486.4546
$ws.Range("K20").Value = 2486.4546
$ws.Range("M20").Value = -2239.4546

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 63832.832
$ws.Range("J141").Value = 62199.4
$ws.Range("L141").Value = 62199.4
$ws.Range("N141").Value = -72559.39999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 59.8
$ws.Range("I6").Value = 66.333336
$ws.Range("K6").Value = 199.000008
$ws.Range("M6").Value = -86.00000800000001
$ws.Range("H11").Value = 767.8333
$ws.Range("I11").Value = 652.5
$ws.Range("K11").Value = 1957.5
$ws.Range("M11").Value = -1817.5
$ws.Range("H107").Value = 291.57144
$ws.Range("J107").Value = 240.16667
$ws.Range("L107").Value = 720.50001
$ws.Range("N107").Value = -4560.50001
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H131").Value = 16289.355
$ws.Range("J131").Value = 17016.094
$ws.Range("L131").Value = 51048.28200000001
$ws.Range("N131").Value = -61128.28200000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 909.25
$ws.Range("I97").Value = 945.9231
$ws.Range("J97").Value = 750.3333
$ws.Range("K97").Value = 945.9231
$ws.Range("L97").Value = 750.3333
$ws.Range("M97").Value = -449.9231
$ws.Range("N97").Value = -1742.3333
$ws.Range("H122").Value = 1499.5
$ws.Range("I122").Value = 1079.2
$ws.Range("K122").Value = 3237.6
$ws.Range("M122").Value = -787.6000000000004
$ws.Range("H132").Value = 1750958.2
$ws.Range("I132").Value = 2138794.2
$ws.Range("K132").Value = 6416382.600000001
$ws.Range("M132").Value = -6413852.600000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 440.2353
$ws.Range("I55").Value = 326.85715
$ws.Range("K55").Value = 326.85715
$ws.Range("M55").Value = -153.85715
$ws.Range("H132").Value = 2188.1133
$ws.Range("I132").Value = 1359.6154
$ws.Range("J132").Value = 2985.926
$ws.Range("K132").Value = 4078.8462
$ws.Range("L132").Value = 8957.778
$ws.Range("M132").Value = -1548.8462
$ws.Range("N132").Value = -14017.778

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H123").Value = 47999
$ws.Range("J123").Value = 47999
$ws.Range("L123").Value = 47999
$ws.Range("N123").Value = -57799
$ws.Range("H136").Value = 12921062
$ws.Range("I136").Value = 20576946
$ws.Range("J136").Value = 1755
$ws.Range("K136").Value = 61730838
$ws.Range("L136").Value = 5265
$ws.Range("M136").Value = -61728288
$ws.Range("N136").Value = -10365
